# Papildināts kļūdu un labojumu excelis
# Adds new rows (23-30) of bug-tracker entries to Sheet1, editing the
# existing row 22's ("22.") follow-up text and inserting 7 fresh rows
# before the previous trailer/legend block (which shifts down by 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert 7 new blank rows right after the current row 23 (i.e. at
#    position 24..30). Excel copies the formatting of the row above
#    (row 23) into the newly inserted rows, which matches the target
#    styling for the new entries.
# ---------------------------------------------------------------------
$ws.Rows("24:30").Insert()

# ---------------------------------------------------------------------
# 2. Fix up row 23 (previously "Nestrādā poga pievienot", dated as a
#    real date serial). It becomes a text date and a reworded comment.
# ---------------------------------------------------------------------
$ws.Range("B23").Value = "08.12.2015."
$ws.Range("C23").Value = "newRoom.php"
$ws.Range("D23").Value = "1. Nestrādā poga ""Pievienot"""
$ws.Range("E23").Value = "A"
$ws.Range("F23").Value = "Krists"

# ---------------------------------------------------------------------
# 3. Populate the 4 new data rows (24-27 logically "23-26" in the sheet).
# ---------------------------------------------------------------------
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "12.12.2015."
$ws.Range("C24").Value = "user-pageForUser.php"
$ws.Range("D24").Value = "1. Lauku ""Lietotāja loma"" vajag virs laukiem ""apgūtie kursi, iegūtie diplomi, iegūtie sertifikāti, pasniedzamie kursi"""
$ws.Range("E24").Value = "J"
$ws.Rows(24).RowHeight = 45

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "12.12.2015."
$ws.Range("C25").Value = "profile.php"
$ws.Range("D25").Value = "1. Ielogojoties ar lietotāju, kura tiesības atbilst pasniedzēja līmenim parādās lapa profile.php, kurā attēlots nenostilots lauks ""Mācību GRUPU PLĀNOŠANA"""
$ws.Range("E25").Value = "J"
$ws.Rows(25).RowHeight = 60

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "12.12.2015."
$ws.Range("C26").Value = "index.php"
$ws.Range("D26").Value = "1. Lauks ""Nav ievadīts lietotājvārds un / vai parole!"" jānostilo, lai tas ir zem pogas ""Ienākt sistēmā"""
$ws.Range("E26").Value = "J"
$ws.Rows(26).RowHeight = 45

# ---------------------------------------------------------------------
# 4. Rows 27-30 stay blank apart from their running numbers (26-29).
# ---------------------------------------------------------------------
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29

# ---------------------------------------------------------------------
# 5. Misc cosmetic touch-ups that came along with the edit.
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 18.14

$ws.Range("H31").Select()
